# Commit message: "swapped slides 5 and 6 on the ppt"
#
# The deck's slide list order needs slide 6 (position 6) and slide 5
# (position 5) to trade places. PowerPoint's native "move slide" operation
# (drag-and-drop in Slide Sorter / the Slide.MoveTo method) is the direct
# COM equivalent of this edit: it repositions the <p:sldId> entry for the
# slide within <p:sldIdLst> in presentation.xml without touching the slide
# parts themselves.
$p = $ppt.ActivePresentation

$slide5 = $p.Slides.Item(5)
$slide6 = $p.Slides.Item(6)

# Move what is currently slide 6 up to position 5; this pushes the former
# slide 5 down to position 6, i.e. the two slides swap places.
$slide6.MoveTo(5)
